# Insert a new row at row 16 (shifting existing rows 16-49 down to 17-50)
# and populate it with the new "Mostre" / "La materia della luce" event.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = "Mostre"
$ws.Range("B16").Value = "Modena"
$ws.Range("C16").Value = "via Borelli, 20"
$ws.Range("D16").Value = "2022-06-08T07:15:59+00:00"
$ws.Range("E16").Value = "Mostra delle opere di Mary Palchetti"
$ws.Range("F16").Value = "2022-06-08T07:16:25+00:00"
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = "2022-06-09T07:00:00+00:00"
$ws.Range("I16").Value = "2022-09-09T08:00:00+00:00"
$ws.Range("J16").Value = "https://www.comune.modena.it/api/novita/eventi/2022/la-materia-della-luce/@@images/538da0ef-6a03-45d5-b15c-38ceef277b8d.jpeg"
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = "2022-06-08T07:17:30+00:00"
$ws.Range("M16").Value = "Associazione Oniro"
$ws.Range("N16").Value = " Inaugurazione 9 giugno ore 20.45"
$ws.Range("O16").Value = ""
$ws.Range("P16").Value = ""
$ws.Range("Q16").Value = ""
$ws.Range("R16").Value = ""
$ws.Range("S16").Value = "La materia della luce"
$ws.Range("T16").Value = ""
$ws.Range("U16").Value = ""
$ws.Range("V16").Value = $false
$ws.Range("W16").Value = 41123
$ws.Range("X16").Value = "https://www.comune.modena.it/novita/eventi/2022/la-materia-della-luce"
$ws.Range("Y16").Value = "44,64582"
$ws.Range("Z16").Value = "10,92572"
$ws.Range("AA16").Value = "POINT (10.92572 44.64582)"
